$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# Running all the suites: set the Runmode column (D) to "Y" for every test case row.
$ws.Range("D2:D26").Value = "Y"

# Reflect the resulting selection (whole Runmode column, active cell at the top).
$ws.Range("D2:D26").Select()
